# Natmi following Dr Hou advice
# Update Il15-Il2rg NATMI LR-pair edge weight table with recomputed values
# (Ligand-expressing cells / Receptor-expressing cells bumped 1 -> 3, and all
# downstream expression/specificity/edge-weight figures recalculated accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.542102
$ws.Range("H2").Value = 16.626306
$ws.Range("I2").Value = 0.2361826998234217
$ws.Range("J2").Value = 0.2361826998234217
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 24.365583
$ws.Range("N2").Value = 73.096749
$ws.Range("O2").Value = 0.3097154004536173
$ws.Range("P2").Value = 0.3097154004536173
$ws.Range("Q2").Value = 135.036546275466
$ws.Range("R2").Value = 1215.328916479194
$ws.Range("S2").Value = 0.07314941945602754
$ws.Range("T2").Value = 0.07314941945602753
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.542102
$ws.Range("H3").Value = 16.626306
$ws.Range("I3").Value = 0.2361826998234217
$ws.Range("J3").Value = 0.2361826998234217
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.228158333333333
$ws.Range("N3").Value = 3.684475
$ws.Range("O3").Value = 0.0156113461364245
$ws.Range("P3").Value = 0.0156113461364245
$ws.Range("Q3").Value = 6.806578755483333
$ws.Range("R3").Value = 61.25920879935
$ws.Range("S3").Value = 0.003687129878378683
$ws.Range("T3").Value = 0.003687129878378682
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.542102
$ws.Range("H4").Value = 16.626306
$ws.Range("I4").Value = 0.2361826998234217
$ws.Range("J4").Value = 0.2361826998234217
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 48.48145033333333
$ws.Range("N4").Value = 145.444351
$ws.Range("O4").Value = 0.6162566192058893
$ws.Range("P4").Value = 0.6162566192058893
$ws.Range("Q4").Value = 268.6891428552673
$ws.Range("R4").Value = 2418.202285697406
$ws.Range("S4").Value = 0.1455491521081012
$ws.Range("T4").Value = 0.1455491521081012
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.542102
$ws.Range("H5").Value = 16.626306
$ws.Range("I5").Value = 0.2361826998234217
$ws.Range("J5").Value = 0.2361826998234217
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.595688
$ws.Range("N5").Value = 13.787064
$ws.Range("O5").Value = 0.05841663420406906
$ws.Range("P5").Value = 0.05841663420406905
$ws.Range("Q5").Value = 25.469771656176
$ws.Range("R5").Value = 229.227944905584
$ws.Range("S5").Value = 0.01379699838091427
$ws.Range("T5").Value = 0.01379699838091427
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.525638333333333
$ws.Range("H6").Value = 10.576915
$ws.Range("I6").Value = 0.1502489091986426
$ws.Range("J6").Value = 0.1502489091986426
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 24.365583
$ws.Range("N6").Value = 73.096749
$ws.Range("O6").Value = 0.3097154004536173
$ws.Range("P6").Value = 0.3097154004536173
$ws.Range("Q6").Value = 85.904233438815
$ws.Range("R6").Value = 773.138100949335
$ws.Range("S6").Value = 0.04653440108017677
$ws.Range("T6").Value = 0.04653440108017676
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.525638333333333
$ws.Range("H7").Value = 10.576915
$ws.Range("I7").Value = 0.1502489091986426
$ws.Range("J7").Value = 0.1502489091986426
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.228158333333333
$ws.Range("N7").Value = 3.684475
$ws.Range("O7").Value = 0.0156113461364245
$ws.Range("P7").Value = 0.0156113461364245
$ws.Range("Q7").Value = 4.330042099402777
$ws.Range("R7").Value = 38.970378894625
$ws.Range("S7").Value = 0.002345587728120225
$ws.Range("T7").Value = 0.002345587728120224
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.525638333333333
$ws.Range("H8").Value = 10.576915
$ws.Range("I8").Value = 0.1502489091986426
$ws.Range("J8").Value = 0.1502489091986426
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 48.48145033333333
$ws.Range("N8").Value = 145.444351
$ws.Range("O8").Value = 0.6162566192058893
$ws.Range("P8").Value = 0.6162566192058893
$ws.Range("Q8").Value = 170.9280597507961
$ws.Range("R8").Value = 1538.352537757165
$ws.Range("S8").Value = 0.09259188482212811
$ws.Range("T8").Value = 0.0925918848221281
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.525638333333333
$ws.Range("H9").Value = 10.576915
$ws.Range("I9").Value = 0.1502489091986426
$ws.Range("J9").Value = 0.1502489091986426
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.595688
$ws.Range("N9").Value = 13.787064
$ws.Range("O9").Value = 0.05841663420406906
$ws.Range("P9").Value = 0.05841663420406905
$ws.Range("Q9").Value = 16.20273378084
$ws.Range("R9").Value = 145.82460402756
$ws.Range("S9").Value = 0.008777035568217491
$ws.Range("T9").Value = 0.008777035568217488
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 13.07613666666666
$ws.Range("H10").Value = 39.22841
$ws.Range("I10").Value = 0.5572537750466107
$ws.Range("J10").Value = 0.5572537750466107
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 24.365583
$ws.Range("N10").Value = 73.096749
$ws.Range("O10").Value = 0.3097154004536173
$ws.Range("P10").Value = 0.3097154004536173
$ws.Range("Q10").Value = 318.60769327101
$ws.Range("R10").Value = 2867.46923943909
$ws.Range("S10").Value = 0.172590076092851
$ws.Range("T10").Value = 0.172590076092851
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 13.07613666666666
$ws.Range("H11").Value = 39.22841
$ws.Range("I11").Value = 0.5572537750466107
$ws.Range("J11").Value = 0.5572537750466107
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.228158333333333
$ws.Range("N11").Value = 3.684475
$ws.Range("O11").Value = 0.0156113461364245
$ws.Range("P11").Value = 0.0156113461364245
$ws.Range("Q11").Value = 16.05956621497222
$ws.Range("R11").Value = 144.53609593475
$ws.Range("S11").Value = 0.008699481568081873
$ws.Range("T11").Value = 0.008699481568081873
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 13.07613666666666
$ws.Range("H12").Value = 39.22841
$ws.Range("I12").Value = 0.5572537750466107
$ws.Range("J12").Value = 0.5572537750466107
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 48.48145033333333
$ws.Range("N12").Value = 145.444351
$ws.Range("O12").Value = 0.6162566192058893
$ws.Range("P12").Value = 0.6162566192058893
$ws.Range("Q12").Value = 633.9500703568788
$ws.Range("R12").Value = 5705.550633211908
$ws.Range("S12").Value = 0.3434113274499434
$ws.Range("T12").Value = 0.3434113274499434
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 13.07613666666666
$ws.Range("H13").Value = 39.22841
$ws.Range("I13").Value = 0.5572537750466107
$ws.Range("J13").Value = 0.5572537750466107
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.595688
$ws.Range("N13").Value = 13.787064
$ws.Range("O13").Value = 0.05841663420406906
$ws.Range("P13").Value = 0.05841663420406905
$ws.Range("Q13").Value = 60.09384436535999
$ws.Range("R13").Value = 540.8445992882399
$ws.Range("S13").Value = 0.03255288993573444
$ws.Range("T13").Value = 0.03255288993573444
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.321440333333333
$ws.Range("H14").Value = 3.964321
$ws.Range("I14").Value = 0.05631461593132515
$ws.Range("J14").Value = 0.05631461593132514
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 24.365583
$ws.Range("N14").Value = 73.096749
$ws.Range("O14").Value = 0.3097154004536173
$ws.Range("P14").Value = 0.3097154004536173
$ws.Range("Q14").Value = 32.197664121381
$ws.Range("R14").Value = 289.778977092429
$ws.Range("S14").Value = 0.01744150382456202
$ws.Range("T14").Value = 0.01744150382456202
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.321440333333333
$ws.Range("H15").Value = 3.964321
$ws.Range("I15").Value = 0.05631461593132515
$ws.Range("J15").Value = 0.05631461593132514
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.228158333333333
$ws.Range("N15").Value = 3.684475
$ws.Range("O15").Value = 0.0156113461364245
$ws.Range("P15").Value = 0.0156113461364245
$ws.Range("Q15").Value = 1.622937957386111
$ws.Range("R15").Value = 14.606441616475
$ws.Range("S15").Value = 0.0008791469618437225
$ws.Range("T15").Value = 0.0008791469618437223
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.321440333333333
$ws.Range("H16").Value = 3.964321
$ws.Range("I16").Value = 0.05631461593132515
$ws.Range("J16").Value = 0.05631461593132514
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 48.48145033333333
$ws.Range("N16").Value = 145.444351
$ws.Range("O16").Value = 0.6162566192058893
$ws.Range("P16").Value = 0.6162566192058893
$ws.Range("Q16").Value = 64.06534388896343
$ws.Range("R16").Value = 576.5880950006708
$ws.Range("S16").Value = 0.03470425482571655
$ws.Range("T16").Value = 0.03470425482571654
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 1.321440333333333
$ws.Range("H17").Value = 3.964321
$ws.Range("I17").Value = 0.05631461593132515
$ws.Range("J17").Value = 0.05631461593132514
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.595688
$ws.Range("N17").Value = 13.787064
$ws.Range("O17").Value = 0.05841663420406906
$ws.Range("P17").Value = 0.05841663420406905
$ws.Range("Q17").Value = 6.072927482615999
$ws.Range("R17").Value = 54.65634734354399
$ws.Range("S17").Value = 0.003289710319202861
$ws.Range("T17").Value = 0.003289710319202861
